$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new column K: Answers / Id / Option
$ws.Range("K1").Value = "Answers"
$ws.Range("K2").Value = "Id"
$ws.Range("K3").Value = "Option"

# Update the active selection on the sheet
$ws.Range("E4").Select()
